$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 154, shifting existing rows 154-189 down to 155-190
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new data record
$ws.Range("A154").Value = 9
$ws.Range("B154").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C154").Value = "Metropolitana"
$ws.Range("D154").Value = 44637
$ws.Range("E154").Value = 13
$ws.Range("F154").Value = "Fruta"
$ws.Range("G154").Value = 100101
$ws.Range("H154").Value = "Berries"
$ws.Range("I154").Value = 100101001
$ws.Range("J154").Value = "Arándano (blue)"
$ws.Range("K154").Value = "Sin especificar"
$ws.Range("L154").Value = "Primera"
$ws.Range("M154").Value = 480
$ws.Range("N154").Value = 4000
$ws.Range("O154").Value = 4000
$ws.Range("P154").Value = 4000
$ws.Range("Q154").Value = "$/bandeja 2 kilos"
$ws.Range("R154").Value = "Provincia de Linares"
$ws.Range("S154").Value = 2000
$ws.Range("T154").Value = 2
